# Yaxuan_Li_Inheritance.xlsx - "Completed Verification of Yaxuan Inheritance Work"
#
# The reviewer filled in / corrected the sample data row (row 4) of the
# inheritance-tracking table on Sheet1, then left the selection on F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: sample/verified data -----------------------------------------
# D = Priority (if incomplete)
$ws.Range("D4").Value = "High"
# C = What is it?
$ws.Range("C4").Value = "Image Preprocessing"
# E = Current State
$ws.Range("E4").Value = "Done"
# F = How to Use It
$ws.Range("F4").Value = "Google Colab, use the code I pushed into the gitlab"
# G = Can It Be Replicated by the new team?
$ws.Range("G4").Value = "Yes but image preprocessing part I think has already done."
# H = What's Missing?
$ws.Range("H4").Value = "No"
# I = Assigned to which Continuing Member?
$ws.Range("I4").Value = "No,  already done"
# J = Due Date
$ws.Range("J4").Value = "Already handed over to Muhammad Arslan"
# K = Link
$ws.Range("K4").Value = "Image_processing/Copy_of_resized.ipynb"
# M = Verfied
$ws.Range("M4").Value = "Yes, I have understood the task and it wasn't much complicated"
# N = Status/Degree of Inheritance
$ws.Range("N4").Value = "Yes, the code runs in my system and I have installed all the required dependencies and installations and this part was not long, so it was easy to understand and inherit as well."
# O = Comments
$ws.Range("O4").Value = "It's already completed as said to me by the members leaving the group members this semester but will als try to see if next semester we could come up with more efficient and good way for the Image PreProcessing Part."

# --- Row heights that changed due to the longer wrapped text --------------
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 90

# --- Column widths, widened to fit the new text (inputs pre-compensated for
#     this host's column-width->stored-width quantization, so the saved
#     <col width="..."/> lands as close as possible to the target) ----------
$ws.Columns.Item(1).ColumnWidth = 5.5
$ws.Columns.Item(2).ColumnWidth = 5.5
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Columns.Item(5).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 46.5
$ws.Columns.Item(7).ColumnWidth = 55.333333333333336
$ws.Columns.Item(8).ColumnWidth = 20.666666666666668
$ws.Columns.Item(9).ColumnWidth = 20.666666666666668
$ws.Columns.Item(10).ColumnWidth = 38.666666666666664
$ws.Columns.Item(11).ColumnWidth = 38.166666666666664
$ws.Columns.Item(12).ColumnWidth = 20.666666666666668
$ws.Columns.Item(13).ColumnWidth = 57.333333333333336
$ws.Columns.Item(14).ColumnWidth = 158.66666666666666
$ws.Columns.Item(15).ColumnWidth = 200

# --- Final selection left on F13 -------------------------------------------
$ws.Range("F13").Select()
